$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend header row (row 1) with new columns P and Q, copying O1's format (bold/border/centered style)
$ws.Range("O1").Copy()
$ws.Range("P1:Q1").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(1, 16).Value = 14   # P1
$ws.Cells.Item(1, 17).Value = 15   # Q1

# For data rows 2-25: swap values in columns I/K and M/O (1<->2),
# and add new columns P and Q (value 2, unstyled like other data cells)
for ($r = 2; $r -le 25; $r++) {
    $ws.Cells.Item($r, 9).Value  = 2   # I -> 2 (was 1)
    $ws.Cells.Item($r, 11).Value = 1   # K -> 1 (was 2)
    $ws.Cells.Item($r, 13).Value = 2   # M -> 2 (was 1)
    $ws.Cells.Item($r, 15).Value = 1   # O -> 1 (was 2)
    $ws.Cells.Item($r, 16).Value = 2   # P -> 2 (new)
    $ws.Cells.Item($r, 17).Value = 2   # Q -> 2 (new)
}
